$wb = $excel.ActiveWorkbook

$oldGuid = "e1c70ddc-cc61-4f5d-b78e-9211ff8b10ae"
$newGuid = "879e77d3-ab9b-4ed8-8ae1-6a51bd3c5903"

$oldMdName  = "$oldGuid.md"
$newMdName  = "$newGuid.md"
$oldMdPath  = "e2e\$oldGuid.md"
$newMdPath  = "e2e\$newGuid.md"

$oldZhXlf = "$oldGuid.10f00e945c1fe57982250f64d0ba3dcffcdabe90.zh-cn.xlf"
$newZhXlf = "$newGuid.8044a48629369dd80ffd6be8418e73f80b4fdbe6.zh-cn.xlf"

$oldDeXlf = "$oldGuid.10f00e945c1fe57982250f64d0ba3dcffcdabe90.de-de.xlf"
$newDeXlf = "$newGuid.8044a48629369dd80ffd6be8418e73f80b4fdbe6.de-de.xlf"

$newHoDate   = "2016-08-29 23:01:13"
$newZhHoDate = "2016-08-29 23:01:08"

$hyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ad3fe04da8334412855d9c4ae1ebdfa5189516a/e2e/$oldGuid.md"

function Set-HyperlinkDisplay($ws, $cellAddr, $displayText) {
    $range = $ws.Range($cellAddr)
    $ws.Hyperlinks.Delete()
    $range.Value = $displayText
    $ws.Hyperlinks.Add($range, $hyperlinkAddr, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
    $range.Font.Underline = $true
    $range.Font.Color = 15570276
}

# ---------- Sheet 1: Overview ----------
$ws1 = $wb.Worksheets.Item("Overview")
Set-HyperlinkDisplay $ws1 "B2" $newMdPath
$ws1.Range("A2").Value = $newMdName
$ws1.Range("G2").Value = $newHoDate

# ---------- Sheet 2: zh-cn ----------
$ws2 = $wb.Worksheets.Item("zh-cn")
Set-HyperlinkDisplay $ws2 "A2" $newMdName
$ws2.Range("G2").Value = $newZhXlf
$ws2.Range("H2").Value = $newZhHoDate

# ---------- Sheet 3: de-de ----------
$ws3 = $wb.Worksheets.Item("de-de")
Set-HyperlinkDisplay $ws3 "A2" $newMdName
$ws3.Range("G2").Value = $newDeXlf
$ws3.Range("H2").Value = $newHoDate
